$d = $word.ActiveDocument

# 1) Correction in CV dates: the "Modelling the learning of numbers in
#    children" project originally read "January 2023 - December 2023".
#    Fix the start date to "January 2022".
$rng = $d.Content
$rng.Find.ClearFormatting()
[void]$rng.Find.Execute("uary 2023", $false, $false, $false, $false, $false, $true, 1, $false, "uary 2022", 2)

# 2) Remove the review comments left on the header line (by Vandana Thambi
#    and Tirthankar Mittra) - this also removes the associated
#    commentRangeStart/End and commentReference markers.
while ($d.Comments.Count -gt 0) {
    $d.Comments.Item(1).Delete()
}
